$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update selection on the sheet (active cell C8)
$ws.Range("C8").Select()

# Update values: B6 and C7 from 0.45 -> 0.4
$ws.Range("B6").Value = 0.4
$ws.Range("C7").Value = 0.4
